$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "LLBV3 Header" column header to "LLBV3 Header / Function"
$ws.Range("E1").Value = "LLBV3 Header / Function"

# Fill in the new "Function" annotations in column E for the rows that
# already have an LLBV3 Header/Pin Name in column D.
# (Order matters for reproducing the exact shared-string table layout.)
$ws.Range("E2").Value = "MCP 2515 interrupt on received frames"

$ws.Range("E6").Value = "X3, for power on board"
$ws.Range("E7").Value = "X3, for power on board"
$ws.Range("E8").Value = "E-stop jumper, also X3"

$ws.Range("E18").Value = "X3, for power on board"
$ws.Range("E19").Value = "wheel hall switch header"

$ws.Range("E21").Value = "all SPI devices, SPI header"
$ws.Range("E20").Value = "SPI header (this pin tells the mega to be a slave)"
$ws.Range("E22").Value = "all SPI devices, SPI header"
$ws.Range("E23").Value = "all SPI devices, SPI header"

$ws.Range("E36").Value = "MCP2515 slave selection"
$ws.Range("E37").Value = "DAC slave selection"

$ws.Range("E53").Value = "on-board buzzer"
$ws.Range("E55").Value = "X3, no purpose assgined"
$ws.Range("E57").Value = "X3, for power on board"

# Match the final viewport / selection state left behind by the edit.
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("E58").Select()
